$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("consort")
$ws.PageSetup.PrintArea = ""
Write-Host $wb.Names.Count
